$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D3").Value = "2016-02-17 06:14:20"
$zhcn.Range("G3").Value = "2016-02-17 06:15:03"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D3").Value = "2016-02-17 06:14:30"
$dede.Range("G3").Value = "2016-02-17 06:15:21"
